$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change stim1_color value from "black" to "red"
$ws.Range("C2").Value = "red"

# Change probe1_color value from "white" to "blue"
$ws.Range("F2").Value = "blue"

# Move the active selection to D8 (matches the recorded cursor position in the diff)
$ws.Range("D8").Select()
